$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4406.9165
$ws.Range("I18").Value = 4548.5
$ws.Range("J18").Value = 3699
$ws.Range("K18").Value = 4548.5
$ws.Range("L18").Value = 3699
$ws.Range("M18").Value = -4264.5
$ws.Range("N18").Value = -4267

$ws.Range("H51").Value = 3071.25
$ws.Range("J51").Value = 3242.5
$ws.Range("L51").Value = 3242.5
$ws.Range("N51").Value = -4210.5

$ws.Range("H55").Value = 546
$ws.Range("I55").Value = 274.2
$ws.Range("K55").Value = 274.2
$ws.Range("M55").Value = -60.19999999999999

$ws.Range("H70").Value = 1999.5
$ws.Range("I70").Value = 1999.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5998.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5728.5
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 1999.5
$ws.Range("I73").Value = 1999.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5998.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5062.5
$ws.Range("N73").ClearContents()

$ws.Range("H100").Value = 2633.9
$ws.Range("I100").Value = 1462.8
$ws.Range("K100").Value = 1462.8
$ws.Range("M100").Value = -921.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2065.9092
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226

$ws.Range("H97").Value = 254.41667
$ws.Range("I97").Value = 239.36363
$ws.Range("J97").Value = 420
$ws.Range("K97").Value = 239.36363
$ws.Range("L97").Value = 420
$ws.Range("M97").Value = 256.63637
$ws.Range("N97").Value = -1412

$ws.Range("H106").Value = 21985
$ws.Range("J106").Value = 21985
$ws.Range("L106").Value = 21985
$ws.Range("N106").Value = -24509

$ws.Range("H116").Value = 2065.9092
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -9588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2065.9092
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228

$ws.Range("H86").Value = 19235722
$ws.Range("I86").Value = 6253.6313
$ws.Range("K86").Value = 6253.6313
$ws.Range("M86").Value = -5130.6313

$ws.Range("H89").Value = 19235722
$ws.Range("I89").Value = 6253.6313
$ws.Range("K89").Value = 31268.1565
$ws.Range("M89").Value = -25652.1565

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 6981.2
$ws.Range("I107").Value = 2699.75
$ws.Range("K107").Value = 2699.75
$ws.Range("M107").Value = -779.75

$ws.Range("H134").Value = 963454.4399999999
$ws.Range("I134").Value = 2084.348
$ws.Range("K134").Value = 6253.044
$ws.Range("M134").Value = -3718.044

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1406.1111
$ws.Range("I22").Value = 531.8
$ws.Range("J22").Value = 2499
$ws.Range("K22").Value = 531.8
$ws.Range("L22").Value = 2499
$ws.Range("M22").Value = -181.8
$ws.Range("N22").Value = -3199

$ws.Range("H28").Value = 25547.666
$ws.Range("J28").Value = 25547.666
$ws.Range("L28").Value = 25547.666
$ws.Range("N28").Value = -26037.666

$ws.Range("H43").Value = 21499.5
$ws.Range("J43").Value = 21499.5
$ws.Range("L43").Value = 21499.5
$ws.Range("N43").Value = -21867.5

$ws.Range("H95").Value = 21712.8
$ws.Range("I95").Value = 18567
$ws.Range("J95").Value = 22499.25
$ws.Range("K95").Value = 18567
$ws.Range("L95").Value = 22499.25
$ws.Range("M95").Value = -15821
$ws.Range("N95").Value = -27991.25

$ws.Range("H101").Value = 21499.5
$ws.Range("J101").Value = 21499.5
$ws.Range("L101").Value = 21499.5
$ws.Range("N101").Value = -27989.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125005784
$ws.Range("I80").Value = 285718000
$ws.Range("J80").Value = 7388.222
$ws.Range("K80").Value = 285718000
$ws.Range("L80").Value = 7388.222
$ws.Range("M80").Value = -285717002
$ws.Range("N80").Value = -9384.222

$ws.Range("H83").Value = 125005784
$ws.Range("I83").Value = 285718000
$ws.Range("J83").Value = 7388.222
$ws.Range("K83").Value = 1428590000
$ws.Range("L83").Value = 36941.11
$ws.Range("M83").Value = -1428585008
$ws.Range("N83").Value = -46925.11

$ws.Range("H98").Value = 20002.5
$ws.Range("J98").Value = 20002.5
$ws.Range("L98").Value = 20002.5
$ws.Range("N98").Value = -25992.5

$ws.Range("H126").Value = 3372.4167
$ws.Range("I126").Value = 1560.25
$ws.Range("K126").Value = 4680.75
$ws.Range("M126").Value = -2210.75

$ws.Range("H132").Value = 22729122
$ws.Range("I132").Value = 25001466
$ws.Range("J132").Value = 5675
$ws.Range("K132").Value = 75004398
$ws.Range("L132").Value = 17025
$ws.Range("M132").Value = -75001868
$ws.Range("N132").Value = -22085

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3443.6875
$ws.Range("I7").Value = 3269.1538
$ws.Range("J7").Value = 4200
$ws.Range("K7").Value = 3269.1538
$ws.Range("L7").Value = 4200
$ws.Range("M7").Value = -3157.1538
$ws.Range("N7").Value = -4424

$ws.Range("H22").Value = 731.375
$ws.Range("J22").Value = 837.5
$ws.Range("L22").Value = 837.5
$ws.Range("N22").Value = -1427.5

$ws.Range("H27").Value = 731.375
$ws.Range("J27").Value = 837.5
$ws.Range("L27").Value = 837.5
$ws.Range("N27").Value = -1051.5

$ws.Range("H68").Value = 5915.067
$ws.Range("I68").Value = 937.8333
$ws.Range("K68").Value = 937.8333
$ws.Range("M68").Value = -188.8333

$ws.Range("H71").Value = 5915.067
$ws.Range("I71").Value = 937.8333
$ws.Range("K71").Value = 4689.1665
$ws.Range("M71").Value = -945.1665000000003

$ws.Range("H82").Value = 50001748
$ws.Range("I82").Value = 1068.7273
$ws.Range("K82").Value = 1068.7273
$ws.Range("M82").Value = -707.7273

$ws.Range("H85").Value = 50001748
$ws.Range("I85").Value = 1068.7273
$ws.Range("K85").Value = 1068.7273
$ws.Range("M85").Value = 179.2727

$ws.Range("H97").Value = 41397.25
$ws.Range("J97").Value = 41397.25
$ws.Range("L97").Value = 41397.25
$ws.Range("N97").Value = -43379.25

$ws.Range("H105").Value = 19500
$ws.Range("J105").Value = 19500
$ws.Range("L105").Value = 19500
$ws.Range("N105").Value = -26488

$ws.Range("H126").Value = 3443.6875
$ws.Range("I126").Value = 3269.1538
$ws.Range("J126").Value = 4200
$ws.Range("K126").Value = 9807.4614
$ws.Range("L126").Value = 12600
$ws.Range("M126").Value = -7337.4614
$ws.Range("N126").Value = -17540

$ws.Range("H132").Value = 9780.299999999999
$ws.Range("I132").Value = 8159.8
$ws.Range("K132").Value = 24479.4
$ws.Range("M132").Value = -21949.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H62").Value = 12209127
$ws.Range("I62").Value = 100003540
$ws.Range("J62").Value = 15459.639
$ws.Range("K62").Value = 100003540
$ws.Range("L62").Value = 15459.639
$ws.Range("M62").Value = -100002916
$ws.Range("N62").Value = -16707.639

$ws.Range("H65").Value = 12209127
$ws.Range("I65").Value = 100003540
$ws.Range("J65").Value = 15459.639
$ws.Range("K65").Value = 500017700
$ws.Range("L65").Value = 77298.19499999999
$ws.Range("M65").Value = -500014580
$ws.Range("N65").Value = -83538.19499999999

$ws.Range("H88").Value = 500025000
$ws.Range("J88").Value = 500025000
$ws.Range("L88").Value = 500025000
$ws.Range("N88").Value = -500025812

$ws.Range("H91").Value = 500025000
$ws.Range("J91").Value = 500025000
$ws.Range("L91").Value = 500025000
$ws.Range("N91").Value = -500027808

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H122").Value = 531295.4399999999
$ws.Range("I122").Value = 1152058.1
$ws.Range("J122").Value = 6034.769
$ws.Range("K122").Value = 3456174.3
$ws.Range("L122").Value = 18104.307
$ws.Range("M122").Value = -3453724.3
$ws.Range("N122").Value = -23004.307

$ws.Range("H126").Value = 9805858
$ws.Range("I126").Value = 1975.9231
$ws.Range("K126").Value = 5927.7693
$ws.Range("M126").Value = -3457.7693
